$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "jun"
$ws.Range("C7").Value = 2.5
$ws.Range("D7").Value = 1.5
$ws.Range("E7").Value = 2345
$ws.Range("F7").Value = 1500
$ws.Range("G7").Value = 4500
$ws.Range("H7").Value = 2

$ws.Range("I7").Select()
